# Insert a new data row at row 17 (pushes existing rows 17.. down by one)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new market record
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = "Vega Monumental Concepción"
$ws.Range("C17").Value = "Bíobío"
$ws.Range("D17").Value = 45274
$ws.Range("D17").NumberFormat = $ws.Range("D18").NumberFormat
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 100112031
$ws.Range("G17").Value = "Poroto verde"
$ws.Range("H17").Value = "Magnum"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 90
$ws.Range("K17").Value = 30000
$ws.Range("L17").Value = 30000
$ws.Range("M17").Value = 30000
$ws.Range("N17").Value = "$/malla 25 kilos"
$ws.Range("O17").Value = "Provincia de Limarí"
$ws.Range("P17").Value = 1200
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = "Hortaliza"
